$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("E3").Value = "2016-03-22 18:45:29"
$wsZh.Range("H3").Value = "2016-03-22 18:46:15"

$wsDe.Range("E3").Value = "2016-03-22 18:45:36"
$wsDe.Range("H3").Value = "2016-03-22 18:46:23"
